$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (fal6_cropped): found/correct columns C and D set to [258]
$ws.Range("C3").Value = "[258]"
$ws.Range("D3").Value = "[258]"

# Row 5 (FallingAwayFromCamera): found/correct columns C and D set to [429]
$ws.Range("C5").Value = "[429]"
$ws.Range("D5").Value = "[429]"
